$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update the date in A1 (was 45406 -> 45436, i.e. 24/04/2024 -> 24/05/2024)
$ws.Range("A1").Value = 45436

# Update the prices in D22:D25
$ws.Range("D22").Value = 12264
$ws.Range("D23").Value = 13894
$ws.Range("D24").Value = 18098
$ws.Range("D25").Value = 20048
